# Update "Estado de Cuenta" worksheet: refresh EC database and add part 1
# of new account statement entries (commit: "Actualiza base de datos EC y
# agrega parte 1 de nuevos estado de cuenta").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Header summary figures
# ----------------------------------------------------------------------
$ws.Range("E11").Value = 463181   # VALOR MORA (total)
$ws.Range("C13").Value = 5        # Cant. Trabajadores

# ----------------------------------------------------------------------
# Detail rows (B:Tipo Doc, C:N Doc, D:Nombre, E:Periodo, F:Valor Mora, G:Salario)
# Row 16 (ROBERT BERTEL ALVIZ / 1607) is unchanged.
# ----------------------------------------------------------------------

# Row 17 -> ANGELICA MARGARITA MURILLO VELASQUEZ, periodo 2003
$ws.Range("C17").Value = "1050958536"
$ws.Range("D17").Value = "ANGELICA MARGARITA MURILLO VELASQUEZ"
$ws.Range("E17").Value = "2003"
$ws.Range("F17").Value = 1170
$ws.Range("G17").Value = 877803

# Row 18 -> BERTHA INES VASQUEZ MARTINEZ, periodo 2303
$ws.Range("C18").Value = "1047432978"
$ws.Range("D18").Value = "BERTHA INES VASQUEZ MARTINEZ"
$ws.Range("E18").Value = "2303"
$ws.Range("F18").Value = 8533
$ws.Range("G18").Value = 1600000

# Row 19 -> BEYNER BAYUELO SOLORZANO, periodo 2411
$ws.Range("C19").Value = "1010031275"
$ws.Range("D19").Value = "BEYNER BAYUELO SOLORZANO"
$ws.Range("E19").Value = "2411"
$ws.Range("F19").Value = 42400
$ws.Range("G19").Value = 1060000

# Row 20 -> BEYNER BAYUELO SOLORZANO, periodo 2412
$ws.Range("E20").Value = "2412"

# Row 21 -> BEYNER BAYUELO SOLORZANO, periodo 2501
$ws.Range("E21").Value = "2501"

# Row 22 -> BEYNER BAYUELO SOLORZANO, periodo 2502
$ws.Range("E22").Value = "2502"

# Row 23 -> BEYNER BAYUELO SOLORZANO, periodo 2503
$ws.Range("E23").Value = "2503"

# Row 24 -> BEYNER BAYUELO SOLORZANO, periodo 2504
$ws.Range("E24").Value = "2504"

# Row 25 -> BEYNER BAYUELO SOLORZANO, periodo 2505
$ws.Range("E25").Value = "2505"

# Row 26 -> BEYNER BAYUELO SOLORZANO, periodo 2506
$ws.Range("E26").Value = "2506"

# Row 27 -> BEYNER BAYUELO SOLORZANO, periodo 2507
$ws.Range("E27").Value = "2507"

# Row 28 -> LAURIN JULIETH PADILLA BARBOZA, periodo 2507
$ws.Range("C28").Value = "1001900489"
$ws.Range("D28").Value = "LAURIN JULIETH PADILLA BARBOZA"
$ws.Range("E28").Value = "2507"
$ws.Range("F28").Value = 1898
$ws.Range("G28").Value = 1423500

# Row 29 -> BEYNER BAYUELO SOLORZANO, periodo 2508 (new statement entry)
$ws.Range("C29").Value = "1010031275"
$ws.Range("D29").Value = "BEYNER BAYUELO SOLORZANO"
$ws.Range("E29").Value = "2508"
$ws.Range("F29").Value = 42400
$ws.Range("G29").Value = 1060000
